$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "305.03"
Set-TextValue "E2" "4.01%"
Set-TextValue "D3" "32.36"
Set-TextValue "E3" "5.76%"
Set-TextValue "D4" "5.317"
Set-TextValue "D5" "0.07597"
Set-TextValue "E5" "6.36%"
Set-TextValue "D6" "1.908"
Set-TextValue "E6" "36.05%"
Set-TextValue "D7" "7.900"
Set-TextValue "E7" "4.46%"
Set-TextValue "D8" "3.877"
Set-TextValue "E8" "7.32%"
Set-TextValue "D9" "0.9321"
Set-TextValue "E9" "1.94%"
Set-TextValue "D10" "0.1715"
Set-TextValue "E10" "4.30%"
Set-TextValue "D11" "0.07982"
Set-TextValue "E11" "2.41%"
Set-TextValue "D12" "0.08104"
Set-TextValue "E12" "4.22%"
Set-TextValue "D13" "0.03058"
Set-TextValue "E13" "3.66%"
Set-TextValue "D14" "0.09949"
Set-TextValue "E14" "10.47%"
Set-TextValue "D15" "0.001501"
Set-TextValue "E15" "-4.75%"
Set-TextValue "D16" "0.04602"
Set-TextValue "E16" "1.41%"
Set-TextValue "D17" "0.006215"
Set-TextValue "E17" "0.88%"
Set-TextValue "D18" "3.445"
Set-TextValue "E18" "-1.15%"
Set-TextValue "D19" "2.232"
Set-TextValue "E19" "-0.52%"
Set-TextValue "D20" "0.3300"
Set-TextValue "E20" "1.44%"
Set-TextValue "D21" "0.1345"
Set-TextValue "E21" "-0.08%"
Set-TextValue "D22" "4.559"
Set-TextValue "E22" "9.95%"
Set-TextValue "D23" "0.1617"
Set-TextValue "E23" "1.74%"
Set-TextValue "D24" "0.001214"
Set-TextValue "E24" "0.51%"
Set-TextValue "D25" "0.004499"
Set-TextValue "E25" "6.17%"
Set-TextValue "E26" "19.75%"
Set-TextValue "D27" "0.0001782"
Set-TextValue "E27" "5.53%"
Set-TextValue "D39" "0.01757"
Set-TextValue "E39" "2,574.96%"
Set-TextValue "D40" "0.04572"
Set-TextValue "E40" "3.58%"
Set-TextValue "D41" "0.007002"
Set-TextValue "E41" "-0.39%"
Set-TextValue "D42" "0.1364"
Set-TextValue "E42" "7.07%"
Set-TextValue "D43" "0.01386"
Set-TextValue "E43" "5.38%"
Set-TextValue "E44" "-5.41%"
Set-TextValue "D45" "0.00006153"
Set-TextValue "E45" "5.11%"
Set-TextValue "D46" "1.896"
Set-TextValue "E46" "0.92%"
Set-TextValue "D47" "0.01222"
Set-TextValue "E47" "-5.95%"
